$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap dates in C10 and C11 (row order was previously 46003/45998, now 45998/46003)
$ws.Range("C10").Value = 45998
$ws.Range("C11").Value = 46003

# Delete row 38 (duplicate "ABC" entry) entirely, shifting dimension to A1:E37
$ws.Rows.Item(38).Delete()
